# Populate previously-empty cells in the phpdocx graph_excel.xlsx template
# with sample/demo data, per the commit "change folder name docx file name
# in template folder and burning issues". All cell styles are left intact;
# only values/text are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shareholding Pattern (rows 4-8): years header + percentages ---
$ws.Range("C4").Value = 2015
$ws.Range("D4").Value = 2014
$ws.Range("E4").Value = 2013
$ws.Range("F4").Value = 2012

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

$ws.Range("C6").Value = 37.96
$ws.Range("D6").Value = 42.1
$ws.Range("E6").Value = 40.52
$ws.Range("F6").Value = 39.02

$ws.Range("C7").Value = 15.1
$ws.Range("D7").Value = -42.1
$ws.Range("E7").Value = 17.51
$ws.Range("F7").Value = 16.57

$ws.Range("C8").Value = 46.94
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 41.97
$ws.Range("F8").Value = 44.41

# --- Liable to retire by rotation counts (rows 16-18) ---
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 0
$ws.Range("C18").Value = 8

# --- ID / NID split (rows 28-29) ---
$ws.Range("C28").Value = 0.8
$ws.Range("D28").Value = 0.2

$ws.Range("C29").Value = 0.8
$ws.Range("D29").Value = 0.2

# --- Executive Compensation vs Shareholder's Value (rows 38-42) ---
$ws.Range("B38").Value = 2011
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 125.66

$ws.Range("B39").Value = 2012
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 115.26

$ws.Range("B40").Value = 2013
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 119.33

$ws.Range("B41").Value = 2014
$ws.Range("C41").Value = 0.83
$ws.Range("D41").Value = 138.61

$ws.Range("B42").Value = 2015
$ws.Range("C42").Value = 6.08
$ws.Range("D42").Value = 186.28

# --- Variation in Director's Remuneration (rows 50-51) ---
$ws.Range("C50").Value = "NA"
$ws.Range("D50").Value = "NA"

$ws.Range("C51").Value = 5.32
$ws.Range("D51").Value = 0.87
